$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H86").Value = 9618.05
$ws.Range("I86").Value = 6688.6665
$ws.Range("J86").Value = 12014.818
$ws.Range("K86").Value = 6688.6665
$ws.Range("L86").Value = 12014.818
$ws.Range("M86").Value = -5565.6665
$ws.Range("N86").Value = -14260.818
$ws.Range("H88").Value = 7022990.5
$ws.Range("I88").Value = 1933.3334
$ws.Range("J88").Value = 7938780.5
$ws.Range("K88").Value = 1933.3334
$ws.Range("L88").Value = 7938780.5
$ws.Range("M88").Value = -1527.3334
$ws.Range("N88").Value = -7939592.5
$ws.Range("H89").Value = 9618.05
$ws.Range("I89").Value = 6688.6665
$ws.Range("J89").Value = 12014.818
$ws.Range("K89").Value = 33443.3325
$ws.Range("L89").Value = 60074.09
$ws.Range("M89").Value = -27827.3325
$ws.Range("N89").Value = -71306.09
$ws.Range("H91").Value = 7022990.5
$ws.Range("I91").Value = 1933.3334
$ws.Range("J91").Value = 7938780.5
$ws.Range("K91").Value = 1933.3334
$ws.Range("L91").Value = 7938780.5
$ws.Range("M91").Value = -529.3334
$ws.Range("N91").Value = -7941588.5

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H28").Value = 1271
$ws.Range("I28").Value = 1271
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1271
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -1079
$ws.Range("H38").Value = 1000
$ws.Range("I38").Value = 1000
$ws.Range("K38").Value = 1000
$ws.Range("M38").Value = -533
$ws.Range("H61").Value = 7814651
$ws.Range("I61").Value = 8930316
$ws.Range("K61").Value = 8930316
$ws.Range("M61").Value = -8930104
$ws.Range("H88").Value = 2901.95
$ws.Range("I88").Value = 2933.3333
$ws.Range("J88").Value = 2876.2727
$ws.Range("K88").Value = 2933.3333
$ws.Range("L88").Value = 2876.2727
$ws.Range("M88").Value = -2527.3333
$ws.Range("N88").Value = -3688.2727
$ws.Range("H91").Value = 2901.95
$ws.Range("I91").Value = 2933.3333
$ws.Range("J91").Value = 2876.2727
$ws.Range("K91").Value = 2933.3333
$ws.Range("L91").Value = 2876.2727
$ws.Range("M91").Value = -1529.3333
$ws.Range("N91").Value = -5684.2727
$ws.Range("H99").Value = 1271
$ws.Range("I99").Value = 1271
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1271
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 1724
$ws.Range("H136").Value = 7814651
$ws.Range("I136").Value = 8930316
$ws.Range("K136").Value = 26790948
$ws.Range("M136").Value = -26788398

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H44").Value = 16000
$ws.Range("J44").Value = 16000
$ws.Range("L44").Value = 16000
$ws.Range("N44").Value = -16994
$ws.Range("H86").Value = 1369950.8
$ws.Range("I86").Value = 1907.75
$ws.Range("J86").Value = 2585989
$ws.Range("K86").Value = 1907.75
$ws.Range("L86").Value = 2585989
$ws.Range("M86").Value = -784.75
$ws.Range("N86").Value = -2588235
$ws.Range("H89").Value = 1369950.8
$ws.Range("I89").Value = 1907.75
$ws.Range("J89").Value = 2585989
$ws.Range("K89").Value = 9538.75
$ws.Range("L89").Value = 12929945
$ws.Range("M89").Value = -3922.75
$ws.Range("N89").Value = -12941177

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 1425
$ws.Range("I31").Value = 889.4737
$ws.Range("J31").Value = 2555.5557
$ws.Range("K31").Value = 889.4737
$ws.Range("L31").Value = 2555.5557
$ws.Range("M31").Value = -594.4737
$ws.Range("N31").Value = -3145.5557
$ws.Range("H34").Value = 1425
$ws.Range("I34").Value = 889.4737
$ws.Range("J34").Value = 2555.5557
$ws.Range("K34").Value = 889.4737
$ws.Range("L34").Value = 2555.5557
$ws.Range("M34").Value = -687.4737
$ws.Range("N34").Value = -2959.5557
$ws.Range("H62").Value = 4139.1875
$ws.Range("I62").Value = 2917
$ws.Range("J62").Value = 4872.5
$ws.Range("K62").Value = 2917
$ws.Range("L62").Value = 4872.5
$ws.Range("M62").Value = -2293
$ws.Range("N62").Value = -6120.5
$ws.Range("H65").Value = 4139.1875
$ws.Range("I65").Value = 2917
$ws.Range("J65").Value = 4872.5
$ws.Range("K65").Value = 14585
$ws.Range("L65").Value = 24362.5
$ws.Range("M65").Value = -11465
$ws.Range("N65").Value = -30602.5

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H40").Value = 578.125
$ws.Range("I40").Value = 217.5
$ws.Range("J40").Value = 698.3333
$ws.Range("K40").Value = 870
$ws.Range("L40").Value = 2793.3332
$ws.Range("M40").Value = -801
$ws.Range("N40").Value = -2931.3332
$ws.Range("H69").Value = 13747
$ws.Range("I69").Value = 356
$ws.Range("J69").Value = 17573
$ws.Range("K69").Value = 1068
$ws.Range("L69").Value = 52719
$ws.Range("M69").Value = -257
$ws.Range("N69").Value = -54341
$ws.Range("H72").Value = 13747
$ws.Range("I72").Value = 356
$ws.Range("J72").Value = 17573
$ws.Range("K72").Value = 3204
$ws.Range("L72").Value = 158157
$ws.Range("M72").Value = 852
$ws.Range("N72").Value = -166269
$ws.Range("H80").Value = 4027.2727
$ws.Range("I80").Value = 1600
$ws.Range("J80").Value = 4270
$ws.Range("K80").Value = 4800
$ws.Range("L80").Value = 12810
$ws.Range("M80").Value = -3864
$ws.Range("N80").Value = -14682
$ws.Range("H83").Value = 4027.2727
$ws.Range("I83").Value = 1600
$ws.Range("J83").Value = 4270
$ws.Range("K83").Value = 14400
$ws.Range("L83").Value = 38430
$ws.Range("M83").Value = -9720
$ws.Range("N83").Value = -47790
$ws.Range("H131").Value = 912.63
$ws.Range("J131").Value = 918.57294
$ws.Range("L131").Value = 2755.71882
$ws.Range("N131").Value = -12835.71882
$ws.Range("H137").Value = 71433120
$ws.Range("I137").Value = 125000950
$ws.Range("J137").Value = 9344.333
$ws.Range("K137").Value = 375002850
$ws.Range("L137").Value = 28032.999
$ws.Range("M137").Value = -374997750
$ws.Range("N137").Value = -38232.999

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 5557621
$ws.Range("I80").Value = 2049.2856
$ws.Range("K80").Value = 2049.2856
$ws.Range("M80").Value = -1051.2856
$ws.Range("H83").Value = 5557621
$ws.Range("I83").Value = 2049.2856
$ws.Range("K83").Value = 10246.428
$ws.Range("M83").Value = -5254.428

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 1437.0526
$ws.Range("I82").Value = 1190.4
$ws.Range("J82").Value = 1711.1111
$ws.Range("K82").Value = 1190.4
$ws.Range("L82").Value = 1711.1111
$ws.Range("M82").Value = -829.4000000000001
$ws.Range("N82").Value = -2433.1111
$ws.Range("H85").Value = 1437.0526
$ws.Range("I85").Value = 1190.4
$ws.Range("J85").Value = 1711.1111
$ws.Range("K85").Value = 1190.4
$ws.Range("L85").Value = 1711.1111
$ws.Range("M85").Value = 57.59999999999991
$ws.Range("N85").Value = -4207.1111

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -183
$ws.Range("H81").Value = 1306.9286
$ws.Range("I81").Value = 1293.875
$ws.Range("J81").Value = 1324.3334
$ws.Range("K81").Value = 2587.75
$ws.Range("L81").Value = 2648.6668
$ws.Range("M81").Value = -1526.75
$ws.Range("N81").Value = -4770.6668
$ws.Range("H84").Value = 1306.9286
$ws.Range("I84").Value = 1293.875
$ws.Range("J84").Value = 1324.3334
$ws.Range("K84").Value = 12938.75
$ws.Range("L84").Value = 13243.334
$ws.Range("M84").Value = -7634.75
$ws.Range("N84").Value = -23851.334
